$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.09990611807518
$ws.Range("C2").Value = 9.284315311862509
$ws.Range("D2").Value = 7.61776766466959
$ws.Range("E2").Value = 13.40717858641336
$ws.Range("F2").Value = 41.27186205150668
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("J2").Value = 10.75331571494002
$ws.Range("K2").Value = 10.32231581683892
$ws.Range("L2").Value = 10.64535173952681
$ws.Range("M2").Value = 15.40973373311131
$ws.Range("O2").Value = 32.0521862894204
$ws.Range("B3").Value = 13.93806356262052
$ws.Range("C3").Value = 9.274518568658557
$ws.Range("D3").Value = 7.609107581246826
$ws.Range("E3").Value = 13.43124580468601
$ws.Range("F3").Value = 41.36873766526165
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("J3").Value = 10.77206884555236
$ws.Range("K3").Value = 10.19932314820355
$ws.Range("L3").Value = 10.65313307650609
$ws.Range("M3").Value = 15.39178170252963
$ws.Range("O3").Value = 32.14194491663513
$ws.Range("B4").Value = 13.8402646854957
$ws.Range("C4").Value = 9.268593908434877
$ws.Range("D4").Value = 7.604719147378981
$ws.Range("E4").Value = 13.44733193144042
$ws.Range("F4").Value = 41.43546384143664
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("J4").Value = 10.7842109304918
$ws.Range("K4").Value = 10.12468246242002
$ws.Range("L4").Value = 10.65895800180019
$ws.Range("M4").Value = 15.382666548072
$ws.Range("O4").Value = 32.20221243681267
$ws.Range("B5").Value = 13.80084981782983
$ws.Range("C5").Value = 9.266202898748881
$ws.Range("D5").Value = 7.603166032683142
$ws.Range("E5").Value = 13.45421679950673
$ws.Range("F5").Value = 41.46447546065713
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("J5").Value = 10.78931718096721
$ws.Range("K5").Value = 10.09451901819092
$ws.Range("L5").Value = 10.66159567621972
$ws.Range("M5").Value = 15.3794350171856
$ws.Range("O5").Value = 32.22806769325163
$ws.Range("B6").Value = 13.79433277544058
$ws.Range("C6").Value = 9.265807299251344
$ws.Range("D6").Value = 7.602922390663476
$ws.Range("E6").Value = 13.45537995194554
$ws.Range("F6").Value = 41.46940270549479
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("J6").Value = 10.79017464098473
$ws.Range("K6").Value = 10.08952664268669
$ws.Range("L6").Value = 10.66204962179672
$ws.Range("M6").Value = 15.37892768666654
$ws.Range("O6").Value = 32.23243918072318
$ws.Range("B7").Value = 13.83973128864623
$ws.Range("C7").Value = 9.26856156721049
$ws.Range("D7").Value = 7.604697247102187
$ws.Range("E7").Value = 13.44742344775384
$ws.Range("F7").Value = 41.43584773457064
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("J7").Value = 10.78427915382818
$ws.Range("K7").Value = 10.12427459946072
$ws.Range("L7").Value = 10.65899250474659
$ws.Range("M7").Value = 15.38262100677428
$ws.Range("O7").Value = 32.20255588423225
$ws.Range("B8").Value = 14.04380254746264
$ws.Range("C8").Value = 9.280918929961187
$ws.Range("D8").Value = 7.614589884316723
$ws.Range("E8").Value = 13.41520565839569
$ws.Range("F8").Value = 41.30376042302009
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("J8").Value = 10.75965180429119
$ws.Range("K8").Value = 10.27974553617622
$ws.Range("L8").Value = 10.64781780792291
$ws.Range("M8").Value = 15.40315003391662
$ws.Range("O8").Value = 32.08206471924887
$ws.Range("B9").Value = 14.45453372779277
$ws.Range("C9").Value = 9.305846646431169
$ws.Range("D9").Value = 7.641284395394727
$ws.Range("E9").Value = 13.36238698546819
$ws.Range("F9").Value = 41.10227605411347
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("J9").Value = 10.71631714592678
$ws.Range("K9").Value = 10.59013934373009
$ws.Range("L9").Value = 10.63418378008816
$ws.Range("M9").Value = 15.45838868697597
$ws.Range("O9").Value = 31.88671340466224
$ws.Range("B10").Value = 14.76012070387167
$ws.Range("C10").Value = 9.324554116184418
$ws.Range("D10").Value = 7.665237897009596
$ws.Range("E10").Value = 13.32986429234046
$ws.Range("F10").Value = 40.98939564755577
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("J10").Value = 10.68747441878231
$ws.Range("K10").Value = 10.81961124016486
$ws.Range("L10").Value = 10.62917447272121
$ws.Range("M10").Value = 15.50788840347538
$ws.Range("O10").Value = 31.7681757178458
$ws.Range("B11").Value = 14.89941877386205
$ws.Range("C11").Value = 9.333143451776484
$ws.Range("D11").Value = 7.677052779329159
$ws.Range("E11").Value = 13.31642617573411
$ws.Range("F11").Value = 40.94568764511426
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("J11").Value = 10.6749975167232
$ws.Range("K11").Value = 10.92390706306736
$ws.Range("L11").Value = 10.62797391394296
$ws.Range("M11").Value = 15.53229230729208
$ws.Range("O11").Value = 31.71968200161801
$ws.Range("B12").Value = 14.9521628633672
$ws.Range("C12").Value = 9.336406806537861
$ws.Range("D12").Value = 7.681656525692632
$ws.Range("E12").Value = 13.31153201829332
$ws.Range("F12").Value = 40.93023602211944
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("J12").Value = 10.67036496592567
$ws.Range("K12").Value = 10.96335462651937
$ws.Range("L12").Value = 10.62767350664583
$ws.Range("M12").Value = 15.54179996203021
$ws.Range("O12").Value = 31.70209996649013
$ws.Range("B13").Value = 14.94080442905302
$ws.Range("C13").Value = 9.335703516259377
$ws.Range("D13").Value = 7.680659296254409
$ws.Range("E13").Value = 13.31257741756654
$ws.Range("F13").Value = 40.93351488894617
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("J13").Value = 10.67135857511282
$ws.Range("K13").Value = 10.95486150750071
$ws.Range("L13").Value = 10.62773135935208
$ws.Range("M13").Value = 15.53974054803945
$ws.Range("O13").Value = 31.70585180515844
$ws.Range("B14").Value = 14.90375842401244
$ws.Range("C14").Value = 9.3334117179275
$ws.Range("D14").Value = 7.67742894660041
$ws.Range("E14").Value = 13.31601963403896
$ws.Range("F14").Value = 40.94439438871535
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("J14").Value = 10.67461454878848
$ws.Range("K14").Value = 10.92715357042469
$ws.Range("L14").Value = 10.62794611334273
$ws.Range("M14").Value = 15.53306919653305
$ws.Range("O14").Value = 31.7182198511262
$ws.Range("B15").Value = 14.88106466467552
$ws.Range("C15").Value = 9.332009307701473
$ws.Range("D15").Value = 7.675467084426957
$ws.Range("E15").Value = 13.31815341383838
$ws.Range("F15").Value = 40.95120162154924
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("J15").Value = 10.67662091969794
$ws.Range("K15").Value = 10.91017454924671
$ws.Range("L15").Value = 10.62809771503527
$ws.Range("M15").Value = 15.52901735109188
$ws.Range("O15").Value = 31.7258974283116
$ws.Range("B16").Value = 14.75101902675267
$ws.Range("C16").Value = 9.323994323059425
$ws.Range("D16").Value = 7.664484047383713
$ws.Range("E16").Value = 13.33076975967133
$ws.Range("F16").Value = 40.99240591250368
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("J16").Value = 10.68830273459673
$ws.Range("K16").Value = 10.81279050682954
$ws.Range("L16").Value = 10.62927455842428
$ws.Range("M16").Value = 15.50633110047452
$ws.Range("O16").Value = 31.77145421659274
$ws.Range("B17").Value = 14.6712802142855
$ws.Range("C17").Value = 9.319097230760308
$ws.Range("D17").Value = 7.657979876772974
$ws.Range("E17").Value = 13.33885657781288
$ws.Range("F17").Value = 41.01964129121961
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("J17").Value = 10.69563375890387
$ws.Range("K17").Value = 10.75300071237331
$ws.Range("L17").Value = 10.63027209027342
$ws.Range("M17").Value = 15.49289357644832
$ws.Range("O17").Value = 31.80079301923265
$ws.Range("B18").Value = 14.62544609525603
$ws.Range("C18").Value = 9.316288046841406
$ws.Range("D18").Value = 7.654325433080077
$ws.Range("E18").Value = 13.34363562433625
$ws.Range("F18").Value = 41.03602561859221
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("J18").Value = 10.69991099106664
$ws.Range("K18").Value = 10.71860451304386
$ws.Range("L18").Value = 10.63094738110652
$ws.Range("M18").Value = 15.48534251140232
$ws.Range("O18").Value = 31.81817895069599
$ws.Range("B19").Value = 14.60993393993664
$ws.Range("C19").Value = 9.315338211688355
$ws.Range("D19").Value = 7.65310304172811
$ws.Range("E19").Value = 13.3452756794338
$ws.Range("F19").Value = 41.0416965764795
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("J19").Value = 10.70136961221685
$ws.Range("K19").Value = 10.70695844252211
$ws.Range("L19").Value = 10.63119348662502
$ws.Range("M19").Value = 15.48281653831867
$ws.Range("O19").Value = 31.82415327555136
$ws.Range("B20").Value = 14.67976582634845
$ws.Range("C20").Value = 9.31961776354602
$ws.Range("D20").Value = 7.658663312465537
$ws.Range("E20").Value = 13.33798250711656
$ws.Range("F20").Value = 41.01666759188183
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("J20").Value = 10.69484708769959
$ws.Range("K20").Value = 10.75936637748406
$ws.Range("L20").Value = 10.63015539849032
$ws.Range("M20").Value = 15.49430565186439
$ws.Range("O20").Value = 31.79761695880686
$ws.Range("B21").Value = 14.91464023381177
$ws.Range("C21").Value = 9.334084587687832
$ws.Range("D21").Value = 7.678374277602529
$ws.Range("E21").Value = 13.31500329529573
$ws.Range("F21").Value = 40.94116896692339
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("J21").Value = 10.67365569176192
$ws.Range("K21").Value = 10.935293603971
$ws.Range("L21").Value = 10.6278788562037
$ws.Range("M21").Value = 15.53502154289195
$ws.Range("O21").Value = 31.71456584171189
$ws.Range("B22").Value = 15.06809401538746
$ws.Range("C22").Value = 9.343601913521164
$ws.Range("D22").Value = 7.69201140258759
$ws.Range("E22").Value = 13.30111888647863
$ws.Range("F22").Value = 40.89823602388608
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("J22").Value = 10.66034301097957
$ws.Range("K22").Value = 11.04998317785384
$ws.Range("L22").Value = 10.62728949857715
$ws.Range("M22").Value = 15.56318228160478
$ws.Range("O22").Value = 31.66484249653331
$ws.Range("B23").Value = 14.98621241412572
$ws.Range("C23").Value = 9.338516835961585
$ws.Range("D23").Value = 7.684664743372369
$ws.Range("E23").Value = 13.30842568243912
$ws.Range("F23").Value = 40.92056346205297
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("J23").Value = 10.66739922257523
$ws.Range("K23").Value = 10.98880858575242
$ws.Range("L23").Value = 10.62752211930937
$ws.Range("M23").Value = 15.54801213560904
$ws.Range("O23").Value = 31.6909637499907
$ws.Range("B24").Value = 14.67592945024586
$ws.Range("C24").Value = 9.319382411255289
$ws.Range("D24").Value = 7.658354066600287
$ws.Range("E24").Value = 13.33837727004677
$ws.Range("F24").Value = 41.01800973888056
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("J24").Value = 10.69520254694839
$ws.Range("K24").Value = 10.75648852701019
$ws.Range("L24").Value = 10.63020783770178
$ws.Range("M24").Value = 15.49366670913713
$ws.Range("O24").Value = 31.79905124015711
$ws.Range("B25").Value = 14.34256008490477
$ws.Range("C25").Value = 9.299031775170265
$ws.Range("D25").Value = 7.633292361058448
$ws.Range("E25").Value = 13.3755700567412
$ws.Range("F25").Value = 41.1506146365518
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("J25").Value = 10.72751229928241
$ws.Range("K25").Value = 10.50577972462734
$ws.Range("L25").Value = 10.63698974796024
$ws.Range("M25").Value = 15.44186281151314
$ws.Range("O25").Value = 31.93517540902371
